$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.4008259790159912
$ws.Range("D2").Value = 0.08476179735234268
$ws.Range("E2").Value = 0.1646544839845525
$ws.Range("F2").Value = 3.188501816221475
$ws.Range("G2").Value = 0.002534671248448313
$ws.Range("J2").Value = 0.2847810779158451
$ws.Range("L2").Value = 0.1424449109324222
$ws.Range("N2").Value = 2.873581687696401
$ws.Range("O2").Value = 9.126135068039105
$ws.Range("C3").Value = 0.3977702695211889
$ws.Range("D3").Value = 0.08500181780032889
$ws.Range("E3").Value = 0.1644652704299574
$ws.Range("F3").Value = 3.139123629587615
$ws.Range("G3").Value = 0.00254060604003649
$ws.Range("J3").Value = 0.2833444913860816
$ws.Range("L3").Value = 0.1428503099764242
$ws.Range("N3").Value = 2.562605684679681
$ws.Range("O3").Value = 8.963443070047845
$ws.Range("C4").Value = 0.3960940797316255
$ws.Range("D4").Value = 0.08516658884908956
$ws.Range("E4").Value = 0.1644172740985681
$ws.Range("F4").Value = 3.110605269914089
$ws.Range("G4").Value = 0.002544444145877141
$ws.Range("J4").Value = 0.2825991955901372
$ws.Range("L4").Value = 0.143144373227603
$ws.Range("N4").Value = 2.371325805375818
$ws.Range("O4").Value = 8.868843805375093
$ws.Range("C5").Value = 0.3954612393418557
$ws.Range("D5").Value = 0.08523811795603109
$ws.Range("E5").Value = 0.164414856852666
$ws.Range("F5").Value = 3.099435058787435
$ws.Range("G5").Value = 0.002546057191811728
$ws.Range("J5").Value = 0.2823298260261424
$ws.Range("L5").Value = 0.1432755714428282
$ws.Range("N5").Value = 2.293303068607429
$ws.Range("O5").Value = 8.831620452771119
$ws.Range("C6").Value = 0.3953591876491345
$ws.Range("D6").Value = 0.08525026032782179
$ws.Range("E6").Value = 0.1644154907085138
$ws.Range("F6").Value = 3.097607468376367
$ws.Range("G6").Value = 0.002546328000479764
$ws.Range("J6").Value = 0.282287170362288
$ws.Range("L6").Value = 0.1432980436041973
$ws.Range("N6").Value = 2.280343261403573
$ws.Range("O6").Value = 8.825519532964677
$ws.Range("C7").Value = 0.3960853417838592
$ws.Range("D7").Value = 0.08516753575444369
$ws.Range("E7").Value = 0.1644171720963641
$ws.Range("F7").Value = 3.110452799291096
$ws.Range("G7").Value = 0.002544465701405622
$ws.Range("J7").Value = 0.2825954237807693
$ws.Range("L7").Value = 0.1431460965786577
$ws.Range("N7").Value = 2.370273851395496
$ws.Range("O7").Value = 8.868336433054878
$ws.Range("C8").Value = 0.3997308048376453
$ws.Range("D8").Value = 0.08484095077592357
$ws.Range("E8").Value = 0.1645750895357239
$ws.Range("F8").Value = 3.171101753324677
$ws.Range("G8").Value = 0.002536677379907218
$ws.Range("J8").Value = 0.2842573246249032
$ws.Range("L8").Value = 0.142575329664016
$ws.Range("N8").Value = 2.766433886209825
$ws.Range("O8").Value = 9.068937128403149
$ws.Range("C9").Value = 0.4084711701628123
$ws.Range("D9").Value = 0.08433819011484012
$ws.Range("E9").Value = 0.165426127758483
$ws.Range("F9").Value = 3.304390414962597
$ws.Range("G9").Value = 0.002522936840162652
$ws.Range("J9").Value = 0.2886043149122202
$ws.Range("L9").Value = 0.1418137995036162
$ws.Range("N9").Value = 3.540180268007646
$ws.Range("O9").Value = 9.504568119725718
$ws.Range("C10").Value = 0.4158706283179754
$ws.Range("D10").Value = 0.08405226851675707
$ws.Range("E10").Value = 0.1663822655072309
$ws.Range("F10").Value = 3.411187870280742
$ws.Range("G10").Value = 0.002513764743261909
$ws.Range("J10").Value = 0.2924660023288226
$ws.Range("L10").Value = 0.1414718669563584
$ws.Range("N10").Value = 4.10623028343673
$ws.Range("O10").Value = 9.850776479272156
$ws.Range("C11").Value = 0.4194509569323941
$ws.Range("D11").Value = 0.08394022163370352
$ws.Range("E11").Value = 0.1668893041853643
$ws.Range("F11").Value = 3.461725954367182
$ws.Range("G11").Value = 0.002509790191109113
$ws.Range("J11").Value = 0.2943689193268
$ws.Range("L11").Value = 0.1413634587048875
$ws.Range("N11").Value = 4.363110593465422
$ws.Range("O11").Value = 10.01404308703718
$ws.Range("C12").Value = 0.4208376719516025
$ws.Range("D12").Value = 0.08390037589296284
$ws.Range("E12").Value = 0.167091686474599
$ws.Range("F12").Value = 3.481146546516129
$ws.Range("G12").Value = 0.002508313404761875
$ws.Range("J12").Value = 0.2951106083343973
$ws.Range("L12").Value = 0.141329177173855
$ws.Range("N12").Value = 4.460285735714251
$ws.Range("O12").Value = 10.07670476525391
$ws.Range("C13").Value = 0.4205376410058932
$ws.Range("D13").Value = 0.08390884257879705
$ws.Range("E13").Value = 0.1670476380775909
$ws.Range("F13").Value = 3.476951363661442
$ws.Range("G13").Value = 0.002508630201815418
$ws.Range("J13").Value = 0.2949499329622398
$ws.Range("L13").Value = 0.1413362593277512
$ws.Range("N13").Value = 4.439361943450422
$ws.Range("O13").Value = 10.06317217017818
$ws.Range("C14").Value = 0.4195644224442958
$ws.Range("D14").Value = 0.08393689175557384
$ws.Range("E14").Value = 0.1669057462419268
$ws.Range("F14").Value = 3.463318018023926
$ws.Range("G14").Value = 0.002509668129083709
$ws.Range("J14").Value = 0.29442951533062
$ws.Range("L14").Value = 0.1413605027005822
$ws.Range("N14").Value = 4.371107314139522
$ws.Range("O14").Value = 10.01918150971164
$ws.Range("C15").Value = 0.4189723282736111
$ws.Range("D15").Value = 0.08395440898297934
$ws.Range("E15").Value = 0.1668201852489979
$ws.Range("F15").Value = 3.455004096420254
$ws.Range("G15").Value = 0.002510307568727561
$ws.Range("J15").Value = 0.2941134938963756
$ws.Range("L15").Value = 0.1413762339279678
$ws.Range("N15").Value = 4.329286057409945
$ws.Range("O15").Value = 9.992345025575219
$ws.Range("C16").Value = 0.4156409665114609
$ws.Range("D16").Value = 0.08405995302973324
$ws.Range("E16").Value = 0.1663505808914856
$ws.Range("F16").Value = 3.407924582323488
$ws.Range("G16").Value = 0.002514028457332159
$ws.Range("J16").Value = 0.2923445903929291
$ws.Range("L16").Value = 0.141479899532186
$ws.Range("N16").Value = 4.089429168003562
$ws.Range("O16").Value = 9.840223313780314
$ws.Range("C17").Value = 0.4136522424130078
$ws.Range("D17").Value = 0.08412931105620203
$ws.Range("E17").Value = 0.1660809650594643
$ws.Range("F17").Value = 3.37954486772054
$ws.Range("G17").Value = 0.002516361666746031
$ws.Range("J17").Value = 0.2912969238347642
$ws.Range("L17").Value = 0.1415555624033189
$ws.Range("N17").Value = 3.94211849063862
$ws.Range("O17").Value = 9.748384667770381
$ws.Range("C18").Value = 0.4125285430194197
$ws.Range("D18").Value = 0.08417090054102516
$ws.Range("E18").Value = 0.1659326735586824
$ws.Range("F18").Value = 3.36340555053269
$ws.Range("G18").Value = 0.002517722302412588
$ws.Range("J18").Value = 0.2907080919579386
$ws.Range("L18").Value = 0.1416035190528255
$ws.Range("N18").Value = 3.857331695637754
$ws.Range("O18").Value = 9.696104696330508
$ws.Range("C19").Value = 0.4121515372976887
$ws.Range("D19").Value = 0.08418527364110773
$ws.Range("E19").Value = 0.1658836294876131
$ws.Range("F19").Value = 3.357972595402202
$ws.Range("G19").Value = 0.002518186195808346
$ws.Range("J19").Value = 0.2905110844965932
$ws.Range("L19").Value = 0.1416205187273789
$ws.Range("N19").Value = 3.828614786364199
$ws.Range("O19").Value = 9.678496715856852
$ws.Range("C20").Value = 0.4138618580061575
$ws.Range("D20").Value = 0.08412175223927321
$ws.Range("E20").Value = 0.1661089638839606
$ws.Range("F20").Value = 3.382546886850889
$ws.Range("G20").Value = 0.002516111365016063
$ws.Range("J20").Value = 0.2914070253062278
$ws.Range("L20").Value = 0.1415470487714856
$ws.Range("N20").Value = 3.95780600327754
$ws.Range("O20").Value = 9.758104789121717
$ws.Range("C21").Value = 0.4198494402158133
$ws.Range("D21").Value = 0.08392858296130612
$ws.Range("E21").Value = 0.1669471415662542
$ws.Range("F21").Value = 3.467314769812788
$ws.Range("G21").Value = 0.002509362497973885
$ws.Range("J21").Value = 0.2945818015051032
$ws.Range("L21").Value = 0.1413531981529239
$ws.Range("N21").Value = 4.391158149571083
$ws.Range("O21").Value = 10.03207989338989
$ws.Range("C22").Value = 0.4239429633021246
$ws.Range("D22").Value = 0.08381739382265607
$ws.Range("E22").Value = 0.16755543470126
$ws.Range("F22").Value = 3.524365370009235
$ws.Range("G22").Value = 0.002505116546253269
$ws.Range("J22").Value = 0.2967796970005168
$ws.Range("L22").Value = 0.1412659647146306
$ws.Range("N22").Value = 4.67379181795809
$ws.Range("O22").Value = 10.21601499353619
$ws.Range("C23").Value = 0.4217416392032192
$ws.Range("D23").Value = 0.0838753621102768
$ws.Range("E23").Value = 0.1672252377327119
$ws.Range("F23").Value = 3.493764842240324
$ws.Range("G23").Value = 0.00250736766349258
$ws.Range("J23").Value = 0.2955953609488091
$ws.Range("L23").Value = 0.1413089149031848
$ws.Range("N23").Value = 4.523002190001307
$ws.Range("O23").Value = 10.11739722765583
$ws.Range("C24").Value = 0.4137670296169915
$ws.Range("D24").Value = 0.08412516423872063
$ws.Range("E24").Value = 0.1660962847028742
$ws.Range("F24").Value = 3.381189124191167
$ws.Range("G24").Value = 0.00251622446654582
$ws.Range("J24").Value = 0.2913572064315986
$ws.Range("L24").Value = 0.1415508839010826
$ws.Range("N24").Value = 3.950713976768498
$ws.Range("O24").Value = 9.753708705221015
$ws.Range("C25").Value = 0.405935500248404
$ws.Range("D25").Value = 0.0844595108407411
$ws.Range("E25").Value = 0.1651378580109757
$ws.Range("F25").Value = 3.266783272651026
$ws.Range("G25").Value = 0.002526491120653027
$ws.Range("J25").Value = 0.2873113850367659
$ws.Range("L25").Value = 0.1419815707235621
$ws.Range("N25").Value = 3.331249627311138
$ws.Range("O25").Value = 9.382154833515585
